# Updates the cryptos list (Coin / Link / Price / Volume(1h) columns),
# as produced by the "Updated cryptos list ... with GitHub Actions" workflow run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the cell to stay a literal string (matches the source export,
    # which always writes Price/Volume as text) instead of letting Excel's
    # type inference turn numeric-looking text into a Number.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") "36.318.68"
$ws.Range("E2").Value = "  -0.17%  "

# Row 3
Set-TextValue $ws.Range("D3") "2.017.56"
$ws.Range("E3").Value = "  -0.70%  "

# Row 4
$ws.Range("E4").Value = "  +0.00%  "

# Row 5
Set-TextValue $ws.Range("D5") "250.58"
$ws.Range("E5").Value = "  +2.77%  "

# Row 6
Set-TextValue $ws.Range("D6") "0.644"
$ws.Range("E6").Value = "  -2.00%  "

# Row 7
Set-TextValue $ws.Range("D7") "63.47"
$ws.Range("E7").Value = "  +19.30%  "

# Row 8
$ws.Range("E8").Value = "  +0.02%  "

# Row 9
Set-TextValue $ws.Range("D9") "59.42"
$ws.Range("E9").Value = "  -2.60%  "

# Row 10
$ws.Range("E10").Value = "  +2.96%  "

# Row 11
Set-TextValue $ws.Range("D11") "0.0749"
$ws.Range("E11").Value = "  +1.55%  "

# Row 12
$ws.Range("E12").Value = "  -0.56%  "

# Row 13
Set-TextValue $ws.Range("D13") "0.941"
$ws.Range("E13").Value = "  +0.01%  "

# Row 14
Set-TextValue $ws.Range("D14") "15.00"
$ws.Range("E14").Value = "  +4.67%  "

# Row 15
Set-TextValue $ws.Range("D15") "2.312.89"
$ws.Range("E15").Value = "  -0.69%  "

# Row 16
$ws.Range("E16").Value = "  +2.82%  "

# Row 17
Set-TextValue $ws.Range("D17") "19.61"
$ws.Range("E17").Value = "  +16.81%  "

# Row 18
Set-TextValue $ws.Range("D18") "2.026.24"
$ws.Range("E18").Value = "  -0.31%  "

# Row 19
Set-TextValue $ws.Range("D19") "36.199.14"
$ws.Range("E19").Value = "  -0.11%  "

# Row 20
Set-TextValue $ws.Range("D20") "72.16"
$ws.Range("E20").Value = "  +1.98%  "

# Row 21
Set-TextValue $ws.Range("D21") "0.0₃0859"
$ws.Range("E21").Value = "  +1.61%  "

# Row 22
Set-TextValue $ws.Range("D22") "5.30"
$ws.Range("E22").Value = "  +3.42%  "

# Row 23
Set-TextValue $ws.Range("D23") "234.73"
$ws.Range("E23").Value = "  -0.79%  "

# Row 24
Set-TextValue $ws.Range("D24") "2.71"
$ws.Range("E24").Value = "  +22.97%  "

# Row 25
$ws.Range("E25").Value = "  +0.12%  "

# Row 26
Set-TextValue $ws.Range("D26") "2.30"
$ws.Range("E26").Value = "  -2.56%  "

# Row 27
Set-TextValue $ws.Range("D27") "9.64"
$ws.Range("E27").Value = "  +6.15%  "

# Row 28
Set-TextValue $ws.Range("D28") "166.09"
$ws.Range("E28").Value = "  +1.92%  "

# Row 29
Set-TextValue $ws.Range("D29") "19.69"
$ws.Range("E29").Value = "  +0.09%  "

# Row 30
$ws.Range("B30").Value = "Filecoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws.Range("D30") "5.25"
$ws.Range("E30").Value = "  +7.46%  "

# Row 31
$ws.Range("B31").Value = "Stellar"
$ws.Range("C31").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue $ws.Range("D31") "0.121"
$ws.Range("E31").Value = "  +0.98%  "

# Row 32
$ws.Range("E32").Value = "  +3.18%  "

# Row 33
Set-TextValue $ws.Range("D33") "0.109"
$ws.Range("E33").Value = "  +26.09%  "

# Row 34
Set-TextValue $ws.Range("D34") "0.0607"
$ws.Range("E34").Value = "  +3.02%  "

# Row 35
Set-TextValue $ws.Range("D35") "4.52"
$ws.Range("E35").Value = "  +4.05%  "

# Row 36
Set-TextValue $ws.Range("D36") "2.46"
$ws.Range("E36").Value = "  +12.58%  "

# Row 37
$ws.Range("E37").Value = "  +0.01%  "

# Row 38
$ws.Range("E38").Value = "  +0.40%  "

# Row 39
Set-TextValue $ws.Range("D39") "5.77"
$ws.Range("E39").Value = "  +18.71%  "

# Row 40
Set-TextValue $ws.Range("D40") "0.109"
$ws.Range("E40").Value = "  +22.82%  "

# Row 41
$ws.Range("E41").Value = "  +0.34%  "

# Row 42
$ws.Range("E42").Value = "  +1.87%  "

# Row 43
$ws.Range("E43").Value = "  +2.60%  "

# Row 44
Set-TextValue $ws.Range("D44") "16.94"
$ws.Range("E44").Value = "  +8.15%  "

# Row 45
$ws.Range("E45").Value = "  +3.58%  "

# Row 46
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue $ws.Range("D46") "7.90"
$ws.Range("E46").Value = "  +6.13%  "

# Row 47
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue $ws.Range("D47") "94.67"
$ws.Range("E47").Value = "  +2.52%  "

# Row 48
Set-TextValue $ws.Range("D48") "1.423.43"
$ws.Range("E48").Value = "  +3.48%  "

# Row 49
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D49") "2.38"
$ws.Range("E49").Value = "  +6.67%  "

# Row 50
$ws.Range("B50").Value = "MXToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue $ws.Range("D50") "2.94"
$ws.Range("E50").Value = "  +1.07%  "

# Row 51
Set-TextValue $ws.Range("D51") "47.35"
$ws.Range("E51").Value = "  +5.41%  "

